$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently follows
#    the H1 title ("Play 3 Stars Slot for Free - Review and Bonuses").
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph ("Play 3 Stars Slot for Free -
#    Review and Bonuses") right before the last paragraph (the
#    "Create a feature image..." image-prompt paragraph).
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 3 Stars Slot for Free - Review and Bonuses</w:t></w:r></w:p>
<w:p/>
<w:sectPr/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertPoint.InsertXML($xml)

# InsertXML leaves behind a spare empty paragraph (from the closing
# <w:p/> used to force a paragraph break) - remove it so only the new
# heading paragraph remains before the image-prompt paragraph.
$spacerPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$spacerRange = $d.Range($spacerPara.Range.Start, $spacerPara.Range.End)
$spacerRange.Delete()

# ------------------------------------------------------------------
# 3) Replace the old image-prompt text with the new meta-description
#    copy, keeping the paragraph's existing (italic) run formatting.
# ------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "Create a feature image for*game development company.",
    $false, $false, $true, $false, $false,
    $true, 1, $false,
    "Experience the rich graphics and progressive jackpots of 3 Stars slot. Read our review and claim free spins, multipliers, and Bonus game rewards.",
    2)

Write-Output "done: $found"
